$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "A."
$ws.Range("D2").Font.Size = 11
$ws.Range("D2").Font.Size = 12
